# Employee export: give every row a distinct "David<n>" name instead of
# the single shared "David", stagger DOJ/Salary per row, and duplicate
# the whole 17-row block as a fresh set of 16 more employees (119-134).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2 and 19 (employees 101 and 118) keep all their other values but
# their shared "David" label becomes "David1" everywhere it is used.
$ws.Range("B2").Value = "David1"
$ws.Range("B19").Value = "David1"

# Rows 3-18 (employees 102-117): each gets its own "David<n>" name, the
# DOJ date advances by one day per row, and Salary increases by 1 per row.
for ($i = 0; $i -lt 16; $i++) {
    $r = 3 + $i
    $ws.Range("B$r").Value = "David" + ($i + 2)
    $ws.Range("G$r").Value = 40311 + $i
    $ws.Range("H$r").Value = 50001 + $i
}

# Seed the date formatting for the new rows by copying the style of an
# existing DOJ cell onto the destination block before filling in values.
$ws.Range("G19").Copy($ws.Range("G20:G35"))

# Append 16 brand-new rows (20-35) for employees 119-134, following the
# exact same progression used for rows 3-18 above.
for ($i = 0; $i -lt 16; $i++) {
    $r = 20 + $i
    $ws.Range("A$r").Value = 119 + $i
    $ws.Range("B$r").Value = "David" + ($i + 2)
    $ws.Range("C$r").Value = 9578821821
    $ws.Range("D$r").Value = 30
    $ws.Range("F$r").Value = "Male"
    $ws.Range("G$r").Value = 40311 + $i
    $ws.Range("H$r").Value = 50001 + $i
}
